$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new test-step rows before the existing "click" row (row 4),
# pushing the old rows 4,5,6 down to 6,7,8.
# ------------------------------------------------------------------
$ws.Rows(4).Insert()
$ws.Rows(4).Insert()

# New row 4: assert the [Google 検索] button is displayed
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "assert"
$ws.Range("C4").Value = "ボタン[Google 検索]表示"
$ws.Range("D4").Value = "is[true]"
$ws.Range("F4").Value = "name[btnK#displayed]"

# New row 5: assert the [Google 検索] button is enabled
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "assert"
$ws.Range("C5").Value = "ボタン[Google 検索]活性"
$ws.Range("D5").Value = "is[true]"
$ws.Range("F5").Value = "name[btnK#enabled]"

# Renumber the "No" column for the rows that shifted down
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# ------------------------------------------------------------------
# Update the print-area defined names to cover the new extent
# ------------------------------------------------------------------
$wb.Names.Item("Sheet1!Print_Area").RefersTo = "=Sheet1!`$A`$1:`$E`$8"
$wb.Names.Item("Sheet1!Print_Area_0").RefersTo = "=Sheet1!`$A`$1:`$E`$7"

# ------------------------------------------------------------------
# Restore the selection to the cell that was active after the edit
# ------------------------------------------------------------------
$ws.Range("F5").Select()
